$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24-36 down to 25-37)
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with the weekly price data
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 45016
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112044
$ws.Range("G24").Value = "Perejil"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 1500
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = 1500
$ws.Range("N24").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 1500
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Range("D24").NumberFormat = $ws.Range("D25").NumberFormat
